# smartart-linear-rule.pptx: move the SmartArt diagram graphic frame down
# on the slide (its vertical offset changes from 1407600 EMU to 2847600 EMU,
# i.e. from 110.8346pt to 224.2205pt -- PowerPoint's Shape.Top is expressed
# in points, so convert EMU -> points by dividing by 12700).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

$sh.Top = 2847600 / 12700
